# Automatic update of files.
# Bump the "Förändrad" (changed) date in column C for all data rows
# (rows 2-16) from 45184 (2023-09-15) to 45185 (2023-09-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C16").Value = 45185
